$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates ---
$ws.Range("E2").Value = 23.39000000000022
$ws.Range("H2").Value = [double]"1.552759475000219e-16"
$ws.Range("K2").Value = 48.03862055582928
$ws.Range("L2").Value = "[45.61044506349436, 50.46679604816419]"
$ws.Range("O2").Value = 1.553500271144502
$ws.Range("P2").Value = "[1.5031844728888082, 1.603816069400196]"
$ws.Range("S2").Value = 52.96021052977162
$ws.Range("T2").Value = "[51.282031631150595, 54.638389428392635]"
$ws.Range("W2").Value = 17.60688688688705
$ws.Range("X2").Value = 17.41957957957974
$ws.Range("Y2").Value = 17.79419419419436

# --- Row 3 updates ---
$ws.Range("E3").Value = 22.76000000000012
$ws.Range("H3").Value = [double]"1.552759475000219e-16"
$ws.Range("I3").Value = ""
$ws.Range("K3").Value = 45.90612209995172
$ws.Range("L3").Value = "[42.840851784501574, 48.97139241540186]"
$ws.Range("O3").Value = -0.3396316382259235
$ws.Range("P3").Value = "[-0.4025263860455386, -0.27673689040630833]"
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 51.35588203451256
$ws.Range("T3").Value = "[49.576511164469174, 53.13525290455595]"
$ws.Range("W3").Value = 1.230270270270278
$ws.Range("X3").Value = 1.002442442442451
$ws.Range("Y3").Value = 1.458098098098106
